$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.00000000222304730179701
$ws.Range("D2").Value = 9844.520545567508
$ws.Range("E2").Value = 616238.5361209477
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 626083.0566665174
